# Applies the cryptos list update as described by the commit diff.
# Source data is text (coin name / link / price / % volume columns), so
# numeric-looking price strings are written with a leading quote to force
# Excel to keep them as text instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.167.29"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "1.642.77"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'217.04"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6
$ws.Range("E6").Value = "  +1.68%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.256"
$ws.Range("E8").Value = "  +1.15%  "

# Row 9
$ws.Range("E9").Value = "  +1.05%  "

# Row 10
$ws.Range("E10").Value = "  +1.29%  "

# Row 11
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.39%  "

# Row 12
$ws.Range("D12").Value = "1.873.06"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").Value = "1.648.50"

# Row 14
$ws.Range("E14").Value = "  +0.49%  "

# Row 15
$ws.Range("D15").Value = "'0.545"
$ws.Range("E15").Value = "  +3.32%  "

# Row 16
$ws.Range("D16").Value = "'67.31"
$ws.Range("E16").Value = "  +1.56%  "

# Row 17
$ws.Range("D17").Value = "27.169.46"
$ws.Range("E17").Value = "  +1.18%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.29%  "

# Row 19
$ws.Range("D19").Value = "'218.67"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20
$ws.Range("E20").Value = "  -0.06%  "

# Row 21
$ws.Range("D21").Value = "'6.83"
$ws.Range("E21").Value = "  +3.05%  "

# Row 22
$ws.Range("E22").Value = "  +5.51%  "

# Row 23
$ws.Range("D23").Value = "'4.41"
$ws.Range("E23").Value = "  +0.86%  "

# Row 24
$ws.Range("E24").Value = "  +0.46%  "

# Row 25
$ws.Range("D25").Value = "'147.70"
$ws.Range("E25").Value = "  +1.42%  "

# Row 26
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'7.53"
$ws.Range("E27").Value = "  +2.11%  "

# Row 28
$ws.Range("E28").Value = "  -0.31%  "

# Row 29
$ws.Range("E29").Value = "  -0.41%  "

# Row 30
$ws.Range("D30").Value = "'0.0508"
$ws.Range("E30").Value = "  -0.58%  "

# Row 31
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +0.20%  "

# Row 32
$ws.Range("E32").Value = "  +0.42%  "

# Row 33
$ws.Range("E33").Value = "  +0.85%  "

# Row 34
$ws.Range("E34").Value = "  +1.24%  "

# Row 35
$ws.Range("D35").Value = "1.262.30"
$ws.Range("E35").Value = "  +2.19%  "

# Row 36
$ws.Range("E36").Value = "  +0.78%  "

# Row 37
$ws.Range("D37").Value = "'0.0177"
$ws.Range("E37").Value = "  +2.25%  "

# Row 38
$ws.Range("E38").Value = "  +0.76%  "

# Row 39
$ws.Range("E39").Value = "  +2.35%  "

# Row 40
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = "  +0.17%  "

# Row 42
$ws.Range("E42").Value = "  +6.35%  "

# Row 43
$ws.Range("D43").Value = "'5.29"
$ws.Range("E43").Value = "  -1.46%  "

# Row 44
$ws.Range("D44").Value = "1.783.77"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").Value = "'61.68"

# Row 46
$ws.Range("D46").Value = "'91.79"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +2.70%  "

# Row 49
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0973"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.59"
$ws.Range("E51").Value = "  +0.07%  "
